$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the header row style (bold / bordered / centered, matching the existing
# "Montreal*" header at AU1) across the new header range before writing values, so
# the new headers reuse the same cell style as the rest of row 1.
$ws.Range("AU1").Copy($ws.Range("AV1:CP1"))

# New bunker-price series header codes for columns AV:CP.
$ws.Range("AV1").Value = "MFSPD00"
$ws.Range("AW1").Value = "MFFJD00"
$ws.Range("AX1").Value = "MFJPD00"
$ws.Range("AY1").Value = "BAMFB00"
$ws.Range("AZ1").Value = "MFSKD00"
$ws.Range("BA1").Value = "WKMFA00"
$ws.Range("BB1").Value = "MFHKD00"
$ws.Range("BC1").Value = "MFSHD00"
$ws.Range("BD1").Value = "MFZSD00"
$ws.Range("BE1").Value = "MFDSY00"
$ws.Range("BF1").Value = "MFDMB00"
$ws.Range("BG1").Value = "MFDKW00"
$ws.Range("BH1").Value = "MFDKF00"
$ws.Range("BI1").Value = "MFDMM00"
$ws.Range("BJ1").Value = "MFDCL00"
$ws.Range("BK1").Value = "MFAGD00"
$ws.Range("BL1").Value = "MFDBD00"
$ws.Range("BM1").Value = "MFGBD00"
$ws.Range("BN1").Value = "MFMLD00"
$ws.Range("BO1").Value = "MFPRD00"
$ws.Range("BP1").Value = "MFRDD00"
$ws.Range("BQ1").Value = "MFDAN00"
$ws.Range("BR1").Value = "MFDGT00"
$ws.Range("BS1").Value = "MFDHB00"
$ws.Range("BT1").Value = "MFDIS00"
$ws.Range("BU1").Value = "MFDLP00"
$ws.Range("BV1").Value = "MFDNV00"
$ws.Range("BW1").Value = "MFDPT00"
$ws.Range("BX1").Value = "MFLIS00"
$ws.Range("BY1").Value = "MFLOM00"
$ws.Range("BZ1").Value = "MFHOD00"
$ws.Range("CA1").Value = "MFNYD00"
$ws.Range("CB1").Value = "MFLAD00"
$ws.Range("CC1").Value = "MFNOD00"
$ws.Range("CD1").Value = "MFPAD00"
$ws.Range("CE1").Value = "MFSED00"
$ws.Range("CF1").Value = "MFVAD00"
$ws.Range("CG1").Value = "MFBAD00"
$ws.Range("CH1").Value = "MFCRD00"
$ws.Range("CI1").Value = "MFSAD00"
$ws.Range("CJ1").Value = "AMFVA00"
$ws.Range("CK1").Value = "AMFCA00"
$ws.Range("CL1").Value = "AMFGY00"
$ws.Range("CM1").Value = "AMFLB00"
$ws.Range("CN1").Value = "AMFMT00"
$ws.Range("CO1").Value = "AMFSF00"
$ws.Range("CP1").Value = "AMFMO00"

# Materialize blank placeholder cells in row 2 (AV2:CP2) for the new series so the
# sheet dimension grows to A1:CP2, mirroring the source data (which has no values
# yet for these new series on this date row). Copying from a genuinely untouched,
# far-away cell stamps blank cells without pulling in any style.
$ws.Range("ZZ999").Copy($ws.Range("AV2:CP2"))

